$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Workbook-level view changes: active tab moves from Sheet3 (index 2) back
#    to Sheet1 (index 0).
# ---------------------------------------------------------------------------
$ws1.Select()

# ---------------------------------------------------------------------------
# 2. Sheet1 view: scroll position / selection update.
# ---------------------------------------------------------------------------
$ws1.Application.ActiveWindow.ScrollRow = 71
$ws1.Range("B109").Select()

# ---------------------------------------------------------------------------
# 3. Row 78 gets a custom row height.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(78).RowHeight = 17.25

# ---------------------------------------------------------------------------
# 4. Helper that builds one monster block starting at $row:
#      row+0 : header  (name, columns C..AZ = 1..50)
#      row+1 : HP      (1.1 / 48)
#      row+2 : Damage  (1.1 / 15)
#      row+3 : Status Effect (None)
#      row+4 : Description
#      row+5 : Basic attack (0.4 / 1)
#      row+6 : 2nd attack   (0.3 / 4)
#      row+7 : 3rd attack   (0.3 / 1.5)
# ---------------------------------------------------------------------------
function Add-MonsterBlock($row, $name, $headerStyleCell, $desc, $atk2Name) {
    $ws1.Range($headerStyleCell).Copy()
    $ws1.Range("A$row").PasteSpecial(-4122)
    $ws1.Range("A$row").Value = $name
    for ($i = 1; $i -le 50; $i++) {
        $ws1.Cells.Item($row, $i + 2).Value = $i
    }

    $r1 = $row + 1
    $ws1.Range("A22").Copy()
    $ws1.Range("A$r1").PasteSpecial(-4122)
    $ws1.Range("A$r1").Value = "HP"
    $ws1.Range("B$r1").Value = 1.1
    $ws1.Range("C$r1").Value = 48

    $r2 = $row + 2
    $ws1.Range("A23").Copy()
    $ws1.Range("A$r2").PasteSpecial(-4122)
    $ws1.Range("A$r2").Value = "Damage"
    $ws1.Range("B$r2").Value = 1.1
    $ws1.Range("C$r2").Value = 15

    $r3 = $row + 3
    $ws1.Range("A24:B24").Copy()
    $ws1.Range("A$r3").PasteSpecial(-4122)
    $ws1.Range("A$r3").Value = "Status Effect"
    $ws1.Range("B$r3").Value = "None"

    $r4 = $row + 4
    $ws1.Range("A25:B25").Copy()
    $ws1.Range("A$r4").PasteSpecial(-4122)
    $ws1.Range("A$r4").Value = "Description"
    $ws1.Range("B$r4").Value = $desc

    $r5 = $row + 5
    $ws1.Range("A70:C70").Copy()
    $ws1.Range("A$r5").PasteSpecial(-4122)
    $ws1.Range("A$r5").Value = "Basic attack"
    $ws1.Range("B$r5").Value = 0.4
    $ws1.Range("C$r5").Value = 1

    $r6 = $row + 6
    $ws1.Range("A71:C71").Copy()
    $ws1.Range("A$r6").PasteSpecial(-4122)
    $ws1.Range("A$r6").Value = $atk2Name
    $ws1.Range("B$r6").Value = 0.3
    $ws1.Range("C$r6").Value = 4

    $r7 = $row + 7
    $ws1.Range("A72:C72").Copy()
    $ws1.Range("A$r7").PasteSpecial(-4122)
    $ws1.Range("A$r7").Value = "Heavy Hit"
    $ws1.Range("B$r7").Value = 0.3
    $ws1.Range("C$r7").Value = 1.5
}

# Order matters: new shared-strings must be created in this exact sequence
# (Cave Spider, Low-med hp..., Poisonous Fangs, Boss - Bandit Leader,
#  Cave Elder, Loyal Minion, Gastropoda) to reproduce the target string table.
Add-MonsterBlock 82 "Cave Spider" "A21" "Low-med hp, med-high dmg" "Poisonous Fangs"
Add-MonsterBlock 109 "Boss - Bandit Leader" "A76" "Low-med hp, med-high dmg" "Poisonous Fangs"
Add-MonsterBlock 91 "Cave Elder" "A21" "Low-med hp, med-high dmg" "Poisonous Fangs"
Add-MonsterBlock 100 "Loyal Minion" "A21" "Low-med hp, med-high dmg" "Poisonous Fangs"
Add-MonsterBlock 118 "Gastropoda" "A21" "Low-med hp, med-high dmg" "Poisonous Fangs"

# ---------------------------------------------------------------------------
# 5. Sheet3 loses the tabSelected flag (handled implicitly by selecting
#    Sheet1 above), keep its own selection untouched.
# ---------------------------------------------------------------------------
$ws3.Range("I9").Select()
$ws1.Select()
$ws1.Range("B109").Select()
